$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 2 ("Suivi du chat forestier") ---
$ws.Range("B2").Value = "Suivi du chat forestier"
$ws.Range("D2").Value = "Présent au sud et à l’est de l’Ile-de-France, le chat forestier est protégé depuis 1979. Son abondance et sa répartition sont étudiés dans le cadre du réseau PMCC. Cependant, des études plus spécifiques sont parfois menées, utilisant des pièges à des points définis."
$ws.Range("E2").Value = "Identifier les zones de présence de l'espèce et évaluer le degré d'hybridation avec le chat domestique"
$ws.Range("F2").Value = "Les résultats des observations sont valorisés en cartes d’abondance et de répartition, articles scientifiques et rapports. Evaluation de l’état de conservation dans le cadre de la directive DHFF."
$ws.Range("J2").Value = "12,1,2,3"
$ws.Range("K2").Value = "Relevés toutes les deux semaines pendant 3 mois"
$ws.Range("L2").Value = "Animation nationale:`nPaul hurel`nSuivi scientifique:`nYoann Bressan`nSandrine Ruette`nAnimateur régional:`nCédric Mondy"
$ws.Range("R2").Value = "Animation`nValidation des observations`nAnalyses génétiques`nBase de données"
$ws.Range("S2").Value = "Animation`nValorisation"
$ws.Range("T2").Value = "Prospections`nRecueil de signalement`nSaisie des données"
$ws.Range("U2").Value = "Formation de 3 jours sur les petits et moyens carnivores et l'utilisation de l'outil rezo-pmcc"
$ws.Range("I2").Value = "Secteurs d'étude d'environ 100km²  comprenant des lisières forestières au contact de prairies. `nChaque secteur est suivi par un minimum de six dispositifs"
$ws.Range("M2").Value = "ONF`nARB`nConseils départementaux (ENS)`nRéserve de la Bassée`nIle-de-France Nature`nCPIE Boucles de la Marne"
$ws.Range("V2").Value = "Signalement de collision:`n- Récupération de tissus et envoi pour analyse`nSuivi pièges photographiques:`n- Pose de dispositifs (piège photo + piège à poil + attractif à base de valériane)`n- Relevés des photos et poils éventuels`n- Stérilisation du piège à poil et recharge en attractif"
$ws.Range("W2").Value = "Signalement de collision:`n- kit de prélèvement PMCC (gants, tube Eppendorf, alcool, ciseaux)`n- fiche adaptée`nSuivi pièges photographiques:`n- pièges photographiques`n- brosses métalliques (pièges à poils)`n- attractif à base de valériane`n- gants, pinces à épiler, enveloppes`n- fiche adaptée"
$ws.Range("Y2").Value = "Saisie des observations sur l'application Rezo-PMCC (pour les observations annexes: saisie sur Rezo-PMCC ou Oison en fonction des espèces)`nTransmission des prélèvements à la DRAS pour analyses génétiques`nValidation des observations sur photo sur la base du phénotype"
$ws.Range("Z2").Value = "SINP national (https://openobs.mnhn.fr) ou régional (https://geonature.arb-idf.fr/geonature/)"
$ws.Range("AA2").Value = "texte:Fiche Espèce;lien:https://professionnels.ofb.fr/fr/doc-fiches-especes/chat-forestier-felis-silvestris-silvestris"
$ws.Range("AB2").Value = "texte:Plaquette de présentation de l'étude;lien:https://oai-gem.ofb.fr/exl-php/document-affiche/ofb_recherche_oai/OUVRE_DOC/49974?fic=doc00073302.pdf"
$ws.Range("AD2").Value = "texte:Site Alfresco de l'étude;lien:https://ged.ofb.fr/share/page/site/etude-chat-forestier-idf/dashboard"
$ws.Range("AE2").Value = "texte:Protocole;lien:https://ged.ofb.fr/share/s/sY4zG36QS1aDJ34fKNlrhw"
$ws.Range("AC2").Value = "texte:Vidéo en collaboration avec le MNHN;lien:https://youtu.be/UopppCJfUHA?feature=shared"
$ws.Range("G2").Value = 77.91
$ws.Range("O2").Value = "1/2j"
$ws.Range("P2").Value = "Maitrise"
$ws.Range("Q2").Value = 1
# Row height for the newly-filled, wrapped-text row
$ws.Rows.Item(2).RowHeight = 225

# Update frozen-pane scroll position / active selection to B2
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B2").Select()
